# Insert "π " into the Engineering Requirement text for the second-mount
# orientation row, turning:
#   "... must be able to vary over at least on at least one axis."
# into:
#   "... must be able to vary over at least π on at least one axis."
#
# The target OOXML keeps the original run intact for the leading portion of
# the sentence and introduces two new runs (same rPr/color) for "π " and the
# trailing portion, so we locate the exact insertion point with Find and then
# materialize the run split explicitly instead of relying on a plain
# Find/Replace (which would just rewrite the text of the single existing
# run).

$d = $word.ActiveDocument

$prefix = "Orientations of the second mount (e.g. receiver mount) must be able to vary over at least "
$suffix = "on at least one axis."
$fullSentence = $prefix + $suffix

$rng = $d.Content
$found = $rng.Find.Execute($fullSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target sentence to edit."
}

$sentenceStart = $rng.Start
$sentenceEnd = $rng.End
$splitPoint = $sentenceStart + $prefix.Length

# Insert the new "π " text right between the existing prefix and suffix.
$insertionPoint = $d.Range($splitPoint, $splitPoint)
$insertionPoint.InsertAfter("π ")

# Toggling a character property on just the newly-inserted text forces it to
# live in its own run (matching the target's 3-run structure) instead of
# being silently re-merged into the surrounding text.
$piRange = $d.Range($splitPoint, $splitPoint + 2)
$piRange.Bold = 1
$piRange.Bold = 0

$finalRange = $d.Range($sentenceStart, $sentenceEnd + 2)
Write-Output "Updated sentence: $($finalRange.Text)"
